# Reworked Loading Sf7 Details
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info ---
$ws.Range("L4").Value = "2020-2021"     # School Year
$ws.Range("P4").Value = "RIZAL"         # Section

# --- Legend / summary counts & prepared-by name ---
$ws.Range("M20").NumberFormat = "@"
$ws.Range("M20").Value = "1"            # MALE count (BoSY)
$ws.Range("M22").NumberFormat = "@"
$ws.Range("M22").Value = "1"            # TOTAL count (BoSY)
$ws.Range("O20").Value = "PHIL REY ESTRELLAÑ PADEROGAO "  # Prepared by

# --- Row 7: replace with former row 8's student data + placeholder columns ---
$ws.Range("A7").Value = "123456789023"
$ws.Range("B7").Value = "Rizal, Jose "
$ws.Range("C7").Value = "M"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2006-05-13"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "14"
$ws.Range("F7").Value = "LOCATION"
$ws.Range("G7").Value = "LANGUAGE"
$ws.Range("H7").Value = "INDIGENOUS PEOPLE"
$ws.Range("I7").Value = "RELIGION"
$ws.Range("J7").Value = "HOUSE NUM"
$ws.Range("K7").Value = "BRGY"
$ws.Range("L7").Value = "MUNICIPAL"
$ws.Range("M7").Value = "PROVINCE"
$ws.Range("N7").Value = "FATHER NAME"
$ws.Range("O7").Value = "MOTHER NAME"
$ws.Range("P7").Value = "GUARDIAN NAME"
$ws.Range("Q7").Value = "RELATIONSHIP"
$ws.Range("R7").Value = "CONTACT"
$ws.Range("S7").Value = "Yolo 3_9"

# --- Row 8: clear entirely (now blank) ---
$ws.Range("A8:S8").ClearContents()
